$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.100.17'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +2.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.466.15'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.11%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.98'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.77'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +3.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.465.67'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +2.10%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.69'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.79%  '

$ws.Range("E11").Value = '  +2.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.406'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +5.49%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.059.71'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.20%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.75'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +6.34%  '

$ws.Range("E15").Value = '  +2.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.469.84'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.58%  '

$ws.Range("E17").Value = '  +1.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.082.01'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.34%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.51%  '

$ws.Range("E20").Value = '  +5.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.31'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +2.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.22'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.561'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.87'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.608.67'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +2.22%  '

$ws.Range("E27").Value = '  +2.27%  '

$ws.Range("E28").Value = '  -1.81%  '

$ws.Range("E29").Value = '  +3.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.996'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.19'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.76%  '

$ws.Range("E32").Value = '  -0.70%  '

$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("E34").Value = '  -0.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.64'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.04%  '

$ws.Range("E36").Value = '  +3.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.10'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.69%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '31.91'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +15.90%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '170.24'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.08%  '

$ws.Range("E40").Value = '  +6.35%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.503.66'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0764'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +1.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.798'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +2.00%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.37'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("E45").Value = '  +5.19%  '

$ws.Range("E46").Value = '  +3.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.42'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.623.86'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +5.96%  '

$ws.Range("E49").Value = '  +12.83%  '

$ws.Range("E50").Value = '  +1.50%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.78'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.49%  '
